# "3d prints. Main screen amended. Last bar issue fixed"
#
# Typography sheet: bump the pixel size of the "Large" and "Huge" fonts.
# Translation sheet: re-point several rows from the "Default" typography
# onto "Tiny" (main-screen labels now rendered in the small/tiny font),
# and append three new translation rows (53-55) for the new 3D-print
# related Label/Unit/Value text ids.

$wb = $excel.ActiveWorkbook

$typography = $wb.Worksheets.Item("Typography")
$translation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: font sizes ---------------------------------------
# Row 5 = "Large" font, row 6 = "Huge" font (column D = Size)
$typography.Range("D5").Value = 66
$typography.Range("D6").Value = 90

# --- Translation sheet: switch several rows from Default to Tiny -------
$rowsToTiny = @(7, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20)
foreach ($r in $rowsToTiny) {
    $translation.Range("C$r").Value = "Tiny"
}

# --- Translation sheet: append new rows 53-55 ---------------------------
$translation.Range("B53").Value = "SingleUseId53"
$translation.Range("C53").Value = "Tiny"
$translation.Range("D53").Value = "Right"
$translation.Range("E53").Value = "LTR"
$translation.Range("F53").Value = "Label"

$translation.Range("B54").Value = "SingleUseId54"
$translation.Range("C54").Value = "Tiny"
$translation.Range("D54").Value = "Left"
$translation.Range("E54").Value = "LTR"
$translation.Range("F54").Value = "Unit"

$translation.Range("B55").Value = "SingleUseId55"
$translation.Range("C55").Value = "Huge"
$translation.Range("D55").Value = "Right"
$translation.Range("E55").Value = "LTR"
$translation.Range("F55").Value = "<value>"
